$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 23 de Junio de 2020 a las 17:16"

# Countries reranked: Moldavia overtakes Ghana (rows 57-58)
$ws.Range("A57").Value = "Moldavia"
$ws.Range("B57").Value = 14714
$ws.Range("C57").Value = 351
$ws.Range("D57").Value = 8212
$ws.Range("E57").Value = 6012
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 10
$ws.Range("H57").Value = 490

$ws.Range("A58").Value = "Ghana"
$ws.Range("B58").Value = 14568
$ws.Range("C58").Value = 414
$ws.Range("D58").Value = 10907
$ws.Range("E58").Value = 3566
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 10
$ws.Range("H58").Value = 95

# Countries reranked: Estado de Palestina overtakes Letonia and Congo (rows 125-127)
$ws.Range("A125").Value = "Estado de Palestina"
$ws.Range("B125").Value = 1157
$ws.Range("C125").Value = 156
$ws.Range("D125").Value = 442
$ws.Range("E125").Value = 712
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 3

$ws.Range("A126").Value = "Letonia"
$ws.Range("B126").Value = 1111
$ws.Range("C126").Value = 0
$ws.Range("D126").Value = 903
$ws.Range("E126").Value = 178
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 30

$ws.Range("A127").Value = "Congo"
$ws.Range("B127").Value = 1087
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 456
$ws.Range("E127").Value = 594
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 37

# Daily case count updates for remaining countries (no reordering)
$ws.Range("B4").Value = 2390273
$ws.Range("C4").Value = 2120
$ws.Range("D4").Value = 1003322
$ws.Range("E4").Value = 1264269
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 72
$ws.Range("H4").Value = 122682

$ws.Range("B7").Value = 449613
$ws.Range("C7").Value = 9163
$ws.Range("D7").Value = 254204
$ws.Range("E7").Value = 181247
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 147
$ws.Range("H7").Value = 14162

$ws.Range("B11").Value = 250767
$ws.Range("C11").Value = 3804
$ws.Range("D11").Value = 210570
$ws.Range("E11").Value = 35692
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 3
$ws.Range("H11").Value = 4505

$ws.Range("B14").Value = 192429
$ws.Range("C14").Value = 310
$ws.Range("D14").Value = 175700
$ws.Range("E14").Value = 7752
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 8
$ws.Range("H14").Value = 8977

$ws.Range("B34").Value = 44931
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 13576
$ws.Range("E34").Value = 30306
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 6
$ws.Range("H34").Value = 1049

$ws.Range("B45").Value = 27936
$ws.Range("C45").Value = 566
$ws.Range("D45").Value = 15551
$ws.Range("E45").Value = 11710
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 6
$ws.Range("H45").Value = 675

$ws.Range("B54").Value = 18231
$ws.Range("C54").Value = 499
$ws.Range("D54").Value = 11220
$ws.Range("E54").Value = 6884
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 127

$ws.Range("B70").Value = 8772
$ws.Range("C70").Value = 21
$ws.Range("D70").Value = 8138
$ws.Range("E70").Value = 386
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 248

$ws.Range("B102").Value = 2318
$ws.Range("C102").Value = 3
$ws.Range("D102").Value = 2123
$ws.Range("E102").Value = 110
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 85

$ws.Range("B110").Value = 1824
$ws.Range("C110").Value = 1
$ws.Range("D110").Value = 1806
$ws.Range("E110").Value = 8
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 10

$ws.Range("B124").Value = 1159
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 1023
$ws.Range("E124").Value = 86
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 50

$ws.Range("B148").Value = 652
$ws.Range("C148").Value = 2
$ws.Range("D148").Value = 270
$ws.Range("E148").Value = 348
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 34
